# update code fix phi coc tu dong tinh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Borrower info
$ws.Range("B7").Value = "ádasdasd"

# Force B8 to stay a text cell (card numbers aren't arithmetic values) while
# keeping its original "style 2" formatting (font/alignment) intact.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "21381723123"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

$ws.Range("B9").Value = "1112000 đồng"

# Signature line (shares the borrower-name string with B7)
$ws.Range("C19").Value = "ádasdasd"

# First borrowed document (row 12)
$ws.Range("B12").Value = "Enzyme Chống Lão Hóa (Tái Bản 2020)"
$ws.Range("C12").Value = "Hiromi Shinya"
$ws.Range("D12").Value = "B-300"

# Second borrowed document (row 13) - new entry
# B13 needs to adopt the same formatting as B12 (the title-style cell)
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = "Street Of Eternal Happiness: Big City Dreams Along A Shanghai Road"
$ws.Range("C13").Value = "Rob Schmitz"
$ws.Range("D13").Value = "B-300"

# Date line
$ws.Range("C16").Value = "TP. Hồ Chí Minh, Ngày 18 tháng 1 năm 2021."
